# Applies the HR report edit described by the commit diff:
#  - Candidate Strengths bullets get a "(high confidence)" qualifier
#    (and some parenthetical detail is trimmed).
#  - Identified Gaps bullets get a coverage qualifier appended
#    (and some parenthetical/detail text is trimmed).
#  - The drafted e-mail body's single QA suggestion is split into two
#    more specific suggestions (joined by a manual line break).
#  - The "Next Steps / QA Suggestions" bullet list gets the same two
#    suggestions as two separate bullets (replacing the old single one).
#
# NOTE: several of the target phrases also appear verbatim in the
# "Detailed Requirement Matches" table; only the ListBullet paragraphs
# under "Candidate Strengths" / "Identified Gaps" are touched, so we
# match on paragraph style + exact text rather than blindly replacing
# every occurrence.

$d = $word.ActiveDocument

function Get-BulletParagraph($doc, $exactText) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Style.NameLocal -ne "List Bullet") { continue }
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $exactText) {
            return $p
        }
    }
    return $null
}

# --- Candidate Strengths bullets -------------------------------------------------

$p = Get-BulletParagraph $d "Strong proficiency in Python, Java, or Go"
$p.Range.Text = "Strong proficiency in Python, Java, or Go (high confidence)."

$p = Get-BulletParagraph $d "Solid understanding of relational databases (PostgreSQL, MySQL) and NoSQL databases (MongoDB, Redis)"
$p.Range.Text = "Solid understanding of relational databases and NoSQL databases (high confidence)."

$p = Get-BulletParagraph $d "Knowledge of containerization technologies (Docker, Kubernetes)"
$p.Range.Text = "Knowledge of containerization technologies (high confidence)."

# --- Identified Gaps bullets -------------------------------------------------

$p = Get-BulletParagraph $d "Experience with version control systems (Git) and CI/CD pipelines"
$p.Range.Text = "Experience with version control systems (Git) and CI/CD pipelines is only moderately covered."

$p = Get-BulletParagraph $d "Bachelor's degree in Computer Science, Engineering, or related field"
$p.Range.Text = "Bachelor's degree in a related field is only moderately covered."

$p = Get-BulletParagraph $d "Knowledge of system design patterns and best practices"
$p.Range.Text = "Knowledge of system design patterns and best practices is the least covered."

# --- Drafted Communication body ----------------------------------------------
# The body is one run containing several <w:t>/<w:br/> children. We locate the
# old sentence by plain-text search (via .IndexOf, not Find.Execute, so
# AutoCorrect/smart-quotes don't mangle the apostrophe we are inserting),
# replace it, then splice in a manual line break followed by the second
# sentence.

$oldSentence = "- Consider asking the candidate about their practical experience with CI/CD pipelines and system design patterns during the interview."
$newSentence1 = "- Clarify the candidate's experience with version control systems and CI/CD pipelines."
$newSentence2 = "- Verify the educational background to ensure it meets the job requirements."

foreach ($p in $d.Paragraphs) {
    $full = $p.Range.Text
    $idx = $full.IndexOf($oldSentence)
    if ($idx -ge 0) {
        $pStart = $p.Range.Start
        $subStart = $pStart + $idx
        $subEnd = $subStart + $oldSentence.Length

        $sub = $d.Range($subStart, $subEnd)
        $sub.Text = $newSentence1

        $afterPos = $subStart + $newSentence1.Length
        $breakRange = $d.Range($afterPos, $afterPos)
        $breakRange.InsertBreak(6)

        $insPos = $afterPos + 1
        $insRange = $d.Range($insPos, $insPos)
        $insRange.InsertAfter($newSentence2)

        break
    }
}

# --- Next Steps / QA Suggestions bullet list ----------------------------------
# Replace the single old bullet (note: no leading "- " here, unlike the
# e-mail body copy above) with the first new sentence, then add a second
# ListBullet paragraph (inherits style) right after it with the second
# sentence.

$oldBullet = "Consider asking the candidate about their practical experience with CI/CD pipelines and system design patterns during the interview."
$newBullet1 = "Clarify the candidate's experience with version control systems and CI/CD pipelines."
$newBullet2 = "Verify the educational background to ensure it meets the job requirements."

$p = Get-BulletParagraph $d $oldBullet
$p.Range.Text = $newBullet1
$p.Range.InsertParagraphAfter()
$p2 = $p.Next()
$p2.Range.Text = $newBullet2

Write-Host "done"
